$wb = $excel.ActiveWorkbook

# --- Environments_OnGoing: append two new rows of data ---
$wsEnv = $wb.Worksheets.Item("Environments_OnGoing")
$wsEnv.Range("A21").Value = "PartnersCommunityOrderRelatedList"
$wsEnv.Range("B21").Value = "/s/order/related/"
$wsEnv.Range("B22").Value = "/csordtelcoa__Orders__r"
$wsEnv.Range("A22").Value = "PartnersCommunityOrderRelatedListView"

# --- PhoneLine: fill in configuration-by-default row + widen column B ---
$wsPhone = $wb.Worksheets.Item("PhoneLine")
$wsPhone.Range("A2").Value = "configurationByDefault"
$wsPhone.Range("B2").Value = "New,NotApplicable,898989,Copper,NotApplicable,NotApplicable,NotApplicable,Classic Telephone Line"
$wsPhone.Columns.Item(2).ColumnWidth = 109.140625

# --- Update selections on each sheet ---
$wsEnv.Range("B31").Select() | Out-Null
$wsPhone.Range("B15").Select() | Out-Null

# --- Make PhoneLine the active sheet/tab ---
$wsPhone.Activate() | Out-Null
